$d = $word.ActiveDocument

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
  <w:pPr>
    <w:jc w:val="center"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="72"/>
      <w:szCs w:val="72"/>
    </w:rPr>
    <w:t>看板方法</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>(Kanban Method)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="72"/>
      <w:szCs w:val="72"/>
    </w:rPr>
    <w:t>基本原則</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>(Foundational Principles)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:ind w:leftChars="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t>從既有的流程開始</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>(Start with existing process)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:ind w:leftChars="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t>同意持續增量、漸進的變化</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>(Agree with pursue incremental, evolutionary change)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:ind w:leftChars="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t>尊重當前的流程、角色、職責和頭銜</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>(Respect the current process, roles , responsibilities and titles)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="a3"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:ind w:leftChars="0"/>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas" w:hint="eastAsia"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t>鼓勵各層級的領導行為</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="56"/>
      <w:szCs w:val="56"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>(Leadership at all level)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas" w:hint="eastAsia"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>「看板」是一種改變管理方針的手段，它不是軟體發展、專案管理的生命週期或是流程，它是在現有軟體發展生命週期或是專案管理方法中引入變化的手段。</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p>
  <w:pPr>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas" w:hint="eastAsia"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>簡單地描述看板方法為</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:b/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>半成品</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:b/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>(WIP)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:b/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>的約束系統</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Consolas" w:eastAsia="微軟正黑體" w:hAnsi="Consolas"/>
      <w:sz w:val="40"/>
      <w:szCs w:val="40"/>
    </w:rPr>
    <w:t>。</w:t>
  </w:r>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($xml)

Write-Output "done"
